# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Refresh the "ESTADO DE CUENTA" worker table (Hoja1!B16:G44) for
# NIT-9004515182: a new worker (HEYDY CAROLINA ORTIZ CABRERA) is added,
# PAOLA ANDREA ARENAS MOGOLLON now covers periods 1905-2105, and
# GLORIA BEATRIZ ARRIETA MEDINA is trimmed down to periods 1906-1908.
#
# Columns: B=Tipo Doc Trabajador, C=N Doc Trabajador, D=Nombre Trabajador,
#          E=Periodo Mora, F=Valor Mora, G=Salario Basico

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

function Set-EcRow {
    param($Row, $TipoDoc, $NumDoc, $Nombre, $Periodo, $ValorMora, $SalarioBasico)
    $ws.Range("B$Row").Value = $TipoDoc
    $ws.Range("C$Row").Value = $NumDoc
    $ws.Range("D$Row").Value = $Nombre
    $ws.Range("E$Row").Value = $Periodo
    $ws.Range("F$Row").Value = $ValorMora
    $ws.Range("G$Row").Value = $SalarioBasico
}

Set-EcRow 16 'CC' '20204308' 'HEYDY CAROLINA ORTIZ CABRERA' '1905' 36000 900000
Set-EcRow 17 'CC' '1030646734' 'PAOLA ANDREA ARENAS MOGOLLON' '1905' 36000 900000
Set-EcRow 18 'CC' '45766004' 'GLORIA BEATRIZ ARRIETA MEDINA' '1906' 33125 877803
Set-EcRow 19 'CC' '1030646734' 'PAOLA ANDREA ARENAS MOGOLLON' '1906' 36000 900000
Set-EcRow 20 'CC' '45766004' 'GLORIA BEATRIZ ARRIETA MEDINA' '1907' 33125 877803
Set-EcRow 21 'CC' '1030646734' 'PAOLA ANDREA ARENAS MOGOLLON' '1907' 36000 900000
Set-EcRow 22 'CC' '45766004' 'GLORIA BEATRIZ ARRIETA MEDINA' '1908' 33125 877803
Set-EcRow 23 'CC' '1030646734' 'PAOLA ANDREA ARENAS MOGOLLON' '1908' 36000 900000
Set-EcRow 24 'CC' '1030646734' 'PAOLA ANDREA ARENAS MOGOLLON' '1909' 36000 900000
Set-EcRow 25 'CC' '1030646734' 'PAOLA ANDREA ARENAS MOGOLLON' '1910' 36000 900000
Set-EcRow 26 'CC' '1030646734' 'PAOLA ANDREA ARENAS MOGOLLON' '1911' 36000 900000
Set-EcRow 27 'CC' '1030646734' 'PAOLA ANDREA ARENAS MOGOLLON' '1912' 36000 900000
Set-EcRow 28 'CC' '1030646734' 'PAOLA ANDREA ARENAS MOGOLLON' '2001' 36000 900000
Set-EcRow 29 'CC' '1030646734' 'PAOLA ANDREA ARENAS MOGOLLON' '2002' 36000 900000
Set-EcRow 30 'CC' '1030646734' 'PAOLA ANDREA ARENAS MOGOLLON' '2003' 36000 900000
Set-EcRow 31 'CC' '1030646734' 'PAOLA ANDREA ARENAS MOGOLLON' '2004' 36000 900000
Set-EcRow 32 'CC' '1030646734' 'PAOLA ANDREA ARENAS MOGOLLON' '2005' 36000 900000
Set-EcRow 33 'CC' '1030646734' 'PAOLA ANDREA ARENAS MOGOLLON' '2006' 36000 900000
Set-EcRow 34 'CC' '1030646734' 'PAOLA ANDREA ARENAS MOGOLLON' '2007' 36000 900000
Set-EcRow 35 'CC' '1030646734' 'PAOLA ANDREA ARENAS MOGOLLON' '2008' 36000 900000
Set-EcRow 36 'CC' '1030646734' 'PAOLA ANDREA ARENAS MOGOLLON' '2009' 36000 900000
Set-EcRow 37 'CC' '1030646734' 'PAOLA ANDREA ARENAS MOGOLLON' '2010' 36000 900000
Set-EcRow 38 'CC' '1030646734' 'PAOLA ANDREA ARENAS MOGOLLON' '2011' 36000 900000
Set-EcRow 39 'CC' '1030646734' 'PAOLA ANDREA ARENAS MOGOLLON' '2012' 36000 900000
Set-EcRow 40 'CC' '1030646734' 'PAOLA ANDREA ARENAS MOGOLLON' '2101' 36000 900000
Set-EcRow 41 'CC' '1030646734' 'PAOLA ANDREA ARENAS MOGOLLON' '2102' 36000 900000
Set-EcRow 42 'CC' '1030646734' 'PAOLA ANDREA ARENAS MOGOLLON' '2103' 36000 900000
Set-EcRow 43 'CC' '1030646734' 'PAOLA ANDREA ARENAS MOGOLLON' '2104' 36000 900000
Set-EcRow 44 'CC' '1030646734' 'PAOLA ANDREA ARENAS MOGOLLON' '2105' 28800 900000
